# DC_Pump_Soft_Starter.xlsx BOM update
# - Update PCB Design doc reference, report date, header label
# - Update diode / MOSFET / relay descriptions & part numbers
# - Split the J1,J2 terminal-block row and R2,R3 resistor row into individual rows
# - Refresh supplier name casing and per-line prices
# - Grow the trailing blank area by two rows

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------------
# 1. Insert the two extra BOM rows first, so everything below lines up with
#    the final layout before we start writing values.
#    Row 10 becomes the new "J2" line (old row 9 "J1, J2" will be split).
#    Row 15 becomes the new "R3" line (old row 13 "R2, R3" will be split).
# ---------------------------------------------------------------------------
$ws.Rows.Item(10).Insert()
$ws.Rows.Item(15).Insert()

# Fix up the formatting of the two newly-inserted (blank) rows so they match
# the alternating banded-row style used throughout the table, instead of the
# default "Normal" style a raw insert leaves behind.
$ws.Range("A8:N8").Copy()
$ws.Range("A10:N10").PasteSpecial(-4122)
$ws.Range("A9:N9").Copy()
$ws.Range("A15:N15").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# ---------------------------------------------------------------------------
# 2. Header / info block
# ---------------------------------------------------------------------------
$ws.Range("C3").Value = "DC_Pump_Soft_Starter.PrjPcb"
$ws.Range("I4").Value = "2.7.2024 г."
$ws.Range("K6").Value = "price"

# ---------------------------------------------------------------------------
# 3. Row 7 - C1 (capacitor): just supplier casing + price
# ---------------------------------------------------------------------------
$ws.Range("I7").Value = "DigiKey"
$ws.Range("K7").Value = "0.31"

# ---------------------------------------------------------------------------
# 4. Row 8 - D1 (diode): new description, new part number, supplier casing,
#    new price
# ---------------------------------------------------------------------------
$ws.Range("B8").Value = " SOD-123 SMD"
$ws.Range("H8").Value = "1N4148W"
$ws.Range("I8").Value = "DigiKey"
$ws.Range("K8").Value = "0.01"

# ---------------------------------------------------------------------------
# 5. Rows 9 & 10 - split "J1, J2" terminal block into two single-designator
#    rows, each with quantity 1 and the refreshed price / supplier casing.
# ---------------------------------------------------------------------------
$ws.Range("D9").Value = 1
$ws.Range("E9").Value = "J1"
$ws.Range("I9").Value = "DigiKey"
$ws.Range("K9").Value = "0.58"

$ws.Range("B10").Value = "TERM BLOCK HDR 2POS 90DEG 3.5MM"
$ws.Range("D10").Value = 1
$ws.Range("E10").Value = "J2"
$ws.Range("G10").Value = "Molex"
$ws.Range("H10").Value = "0395021002"
$ws.Range("I10").Value = "DigiKey"
$ws.Range("J10").Value = "WM7770-ND"
$ws.Range("K10").Value = "0.58"

# ---------------------------------------------------------------------------
# 6. Row 11 - K1 (relay): now has a description / device package
# ---------------------------------------------------------------------------
$ws.Range("B11").Value = "TRKM S-Z L 12VDC"

# ---------------------------------------------------------------------------
# 7. Row 12 - Q1 (MOSFET): updated part description
# ---------------------------------------------------------------------------
$ws.Range("B12").Value = "N-Channel 100V 33A (Tc) 94W (Tc) Through Hole TO-220AB"

# ---------------------------------------------------------------------------
# 8. Rows 14 & 15 - split "R2, R3" resistor row into two single-designator
#    rows, each with quantity 1 and the refreshed price / supplier casing.
# ---------------------------------------------------------------------------
$ws.Range("C14").Value = $null
$ws.Range("D14").Value = 1
$ws.Range("E14").Value = "R2"
$ws.Range("G14").Value = $null
$ws.Range("H14").Value = $null
$ws.Range("I14").Value = $null
$ws.Range("J14").Value = $null
$ws.Range("K14").Value = $null

$ws.Range("B15").Value = "RESHIGHPOWERA 0805 100K 1% 1/2W"
$ws.Range("C15").Value = "0805 (2012 Metric)"
$ws.Range("D15").Value = 1
$ws.Range("E15").Value = "R3"
$ws.Range("G15").Value = "Bourns Inc."
$ws.Range("H15").Value = "CHP0805AFX-1003ELF"
$ws.Range("I15").Value = "DigiKey"
$ws.Range("J15").Value = "118-CHP0805AFX-1003ELFCT-ND"
$ws.Range("K15").Value = "0.06"

# Row 14 keeps the R2 data (qty/descr/package/manufacturer/part/supplier/price)
$ws.Range("B14").Value = "RESHIGHPOWERA 0805 100K 1% 1/2W"
$ws.Range("C14").Value = "0805 (2012 Metric)"
$ws.Range("G14").Value = "Bourns Inc."
$ws.Range("H14").Value = "CHP0805AFX-1003ELF"
$ws.Range("I14").Value = "DigiKey"
$ws.Range("J14").Value = "118-CHP0805AFX-1003ELFCT-ND"
$ws.Range("K14").Value = "0.06"

# ---------------------------------------------------------------------------
# 9. Extend the trailing blank banded area by two rows (40 & 41), matching
#    the style used by the existing blank rows above them.
# ---------------------------------------------------------------------------
$ws.Range("A39:N39").Copy()
$ws.Range("A40:N41").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$wb.Save()
